# Allocation rule updated with 5 and 10 mi rad
# - sq_miles (E) and pop_sq_mile_1mi (G) values are re-expressed as text
#   (same numeric text, stored as shared strings instead of numbers).
# - total_risk (R) and total_risk_resp (S) are recomputed/rounded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  E="20.5752959055157";  G="142.549590220656";  R=50;    S=0.6 },
    @{ Row=3;  E="103.612502034919";  G="45.6605130373706";  R=50;    S=0.6 },
    @{ Row=4;  E="36.0627684203524";  G="166.099283620707";  R=30;    S=0.366666666666667 },
    @{ Row=5;  E="11.777227783969";   G="1405.25430123102";  R=30;    S=0.383333333333333 },
    @{ Row=6;  E="74.6569158825578";  G="53.5114523922272";  R=35;    S=1.45 },
    @{ Row=7;  E="120.261973778333";  G="75.7429777162124";  R=107.5; S=0.525 },
    @{ Row=8;  E="20.8371989696875";  G="79.5212447896903";  R=200;   S=0.6 },
    @{ Row=9;  E="17.2846381134759";  G="419.968295103647";  R=20;    S=0.2 },
    @{ Row=10; E="130.239383628461";  G="40.0646859239259";  R=20;    S=0.2 }
)

foreach ($d in $rows) {
    $eCell = $ws.Cells.Item($d.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $d.E
    $eCell.ClearFormats()

    $gCell = $ws.Cells.Item($d.Row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $d.G
    $gCell.ClearFormats()

    $ws.Cells.Item($d.Row, 18).Value = $d.R
    $ws.Cells.Item($d.Row, 19).Value = $d.S
}
